# ---------------------------------------------------------------------------
# "Book of Itza" review doc:
#   1. Add a new "Meta description: ..." paragraph right after the title.
#   2. Drop the duplicate bold "Play Book of Itza..." paragraph near the end.
#   3. Turn the trailing italic paragraph into an image-generation prompt.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$boldLabelText  = "Play Book of Itza online slot game for free now!"
$italicOldText  = "Read our Book of Itza online slot game review and play for free! Enjoy the amazing Aztec theme and unique gameplay mechanism of this medium-volatility game."

function Find-ParagraphByText($doc, $text, $excludeStarts) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $text) {
            $skip = $false
            foreach ($s in $excludeStarts) {
                if ($p.Range.Start -eq $s) { $skip = $true }
            }
            if (-not $skip) {
                return $p
            }
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Step 1: Insert a new "Meta description" paragraph right after the first
# paragraph ("Play Book of Itza online slot game for free now!" / Heading1).
#
# Build it by duplicating the existing bold "Play Book of Itza..." paragraph
# near the end of the document via Copy/Paste -- this produces a paragraph
# shaped like the target (no paragraph style override, one bold run) -- and
# then rewrite its text and run formatting in place.
# ---------------------------------------------------------------------------

$p1 = $d.Paragraphs.Item(1)
$pBoldTemplate = Find-ParagraphByText $d $boldLabelText @($p1.Range.Start)

$templateRange = $d.Range($pBoldTemplate.Range.Start, $pBoldTemplate.Range.End)
$templateRange.Copy()

$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)
$insertPoint.Paste()

$metaPara = $d.Paragraphs.Item(2)

# Strip bold from the (still single) pasted run before rewriting its text so
# the text we type next doesn't inherit bold.
$metaFull = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$metaFull.Bold = 0

# Replace the run's text with the full sentence (bold label + plain detail
# will be split back out below).
$metaFull2 = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$metaFull2.Text = "Meta description: Read our Book of Itza online slot game review and play for free! Enjoy the amazing Aztec theme and unique gameplay mechanism of this medium-volatility game."

# Re-bold only the "Meta description" label portion.
$labelRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$labelRange.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$labelRange.Bold = 1

# ---------------------------------------------------------------------------
# Step 2: Remove the duplicate "Play Book of Itza online slot game for free
# now!" paragraph (bold) near the end of the document. Re-resolve it by text
# now that offsets have shifted, excluding the title and the freshly
# inserted meta-description paragraph.
# ---------------------------------------------------------------------------

$p1 = $d.Paragraphs.Item(1)
$pDup = Find-ParagraphByText $d $boldLabelText @($p1.Range.Start, $metaPara.Range.Start)
$pDup.Range.Delete()

# ---------------------------------------------------------------------------
# Step 3: Replace the italic closing paragraph's text with the new "Prompt:"
# text (the run stays italic -- only its text content changes). Re-resolve
# it fresh now that the document has changed further.
# ---------------------------------------------------------------------------

$pItalic = Find-ParagraphByText $d $italicOldText @()
$italicRange = $d.Range($pItalic.Range.Start, $pItalic.Range.End - 1)
$italicRange.Text = "Prompt: Create a feature image for `"Book of Itza`" that showcases the happy Maya warrior with glasses in a cartoon style. The image should feature the warrior standing in front of an ancient temple, holding the titular Book of Itza. The temple should have Aztec designs and symbols, and the background should be vibrant and colorful. The warrior should have a big smile on his face and eye-catching details such as feathered headdress and intricate tattoos. The image should convey the excitement and adventure of playing the slot game while highlighting the Aztec theme and the expanding wild feature represented by the Book of Itza."

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
